# Journal_Travail.xlsx — "Ajout de l'action des boutons Nouveau, Ouvrir,
# Sauvegarder et Sauvergarder sous de la menuBar." journal entry.
#
# Adds a new journal row (row 15): date 27/03/2018 (serial 43186),
# an activity description with "Nouveau", "Ouvrir", "Sauvegarder" and
# "Sauvergarder sous" in italics, and 3 hours worked. The Total (C32)
# recalculates automatically via its existing SUM formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15: date / activity / hours -------------------------------------
$ws.Range("A15").Value = 43186

$activityText = "Ajout de l'action des boutons Nouveau, Ouvrir, Sauvegarder et Sauvergarder sous."
$ws.Range("B15").Value = $activityText

$ws.Range("C15").Value = 3

# Row grows to fit the wrapped two-line description, like the other
# multi-line entries above it.
$ws.Range("A15:C15").RowHeight = 30

# --- Rich-text (italic) formatting on the button names --------------------
$cell = $ws.Range("B15")

$segments = @(
    @{ text = "Ajout de l'action des boutons "; italic = -1 },
    @{ text = "Nouveau"; italic = 1 },
    @{ text = ", "; italic = 0 },
    @{ text = "Ouvrir"; italic = 1 },
    @{ text = ", "; italic = 0 },
    @{ text = "Sauvegarder"; italic = 1 },
    @{ text = " et "; italic = 0 },
    @{ text = "Sauvergarder sous"; italic = 1 },
    @{ text = "."; italic = 0 }
)

$pos = 1
foreach ($seg in $segments) {
    $len = $seg.text.Length
    if ($seg.italic -eq 1) {
        $cell.Characters($pos, $len).Font.Italic = $true
    } elseif ($seg.italic -eq 0) {
        $cell.Characters($pos, $len).Font.Italic = $false
    }
    $pos = $pos + $len
}

# --- Selection moves to where the author ended up editing -----------------
$ws.Range("H14").Select()
